$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Soulsweeper", $false, $false, $false, $false, $false, $true, 1, $false, [char]0x2018 + "SoleSweeper", 2)
